# Updated symbol list - applies refreshed Price (D) and Volume(1h) (E) values
# to the cryptos worksheet, matching the upstream GitHub Actions scrape commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to store the value as literal text (matches the inline-string
    # cells already used throughout this sheet) instead of letting Excel auto-convert
    # numeric-looking / percent-looking input into a Number cell, then drop the
    # temporary "@" number format so no stray style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.72"
Set-TextValue $ws.Range("E2") "2.89%"
Set-TextValue $ws.Range("D3") "34.83"
Set-TextValue $ws.Range("E3") "12.34%"
Set-TextValue $ws.Range("D4") "5.162"
Set-TextValue $ws.Range("E4") "4.80%"
Set-TextValue $ws.Range("D5") "0.07745"
Set-TextValue $ws.Range("E5") "4.68%"
Set-TextValue $ws.Range("D6") "2.374"
Set-TextValue $ws.Range("E6") "7.39%"
Set-TextValue $ws.Range("D7") "8.011"
Set-TextValue $ws.Range("E7") "3.60%"
Set-TextValue $ws.Range("D8") "3.948"
Set-TextValue $ws.Range("E8") "5.25%"
Set-TextValue $ws.Range("D9") "0.9296"
Set-TextValue $ws.Range("E9") "1.68%"
Set-TextValue $ws.Range("D10") "0.09788"
Set-TextValue $ws.Range("E10") "9.33%"
Set-TextValue $ws.Range("D11") "0.1796"
Set-TextValue $ws.Range("E11") "5.47%"
Set-TextValue $ws.Range("D12") "0.08634"
Set-TextValue $ws.Range("E12") "4.04%"
Set-TextValue $ws.Range("E13") "6.27%"
Set-TextValue $ws.Range("D14") "0.09898"
Set-TextValue $ws.Range("E14") "-0.85%"
Set-TextValue $ws.Range("D15") "0.001500"
Set-TextValue $ws.Range("E15") "-0.29%"
Set-TextValue $ws.Range("D16") "0.005765"
Set-TextValue $ws.Range("E16") "-1.43%"
Set-TextValue $ws.Range("D17") "3.460"
Set-TextValue $ws.Range("E17") "-1.13%"
Set-TextValue $ws.Range("D18") "2.146"
Set-TextValue $ws.Range("E18") "3.92%"
Set-TextValue $ws.Range("D19") "0.3367"
Set-TextValue $ws.Range("E19") "1.20%"
Set-TextValue $ws.Range("E20") "2.77%"
Set-TextValue $ws.Range("D21") "4.347"
Set-TextValue $ws.Range("E21") "9.17%"
Set-TextValue $ws.Range("D22") "0.2301"
Set-TextValue $ws.Range("E22") "5.20%"
Set-TextValue $ws.Range("D23") "0.04578"
Set-TextValue $ws.Range("E23") "0.43%"
Set-TextValue $ws.Range("D24") "0.001218"
Set-TextValue $ws.Range("E24") "0.36%"
Set-TextValue $ws.Range("D25") "0.004461"
Set-TextValue $ws.Range("E25") "-2.68%"
Set-TextValue $ws.Range("D26") "0.0001301"
Set-TextValue $ws.Range("E26") "-0.06%"
Set-TextValue $ws.Range("E27") "-0.17%"
Set-TextValue $ws.Range("D39") "0.01785"
Set-TextValue $ws.Range("E39") "11.57%"
Set-TextValue $ws.Range("D40") "0.04799"
Set-TextValue $ws.Range("E40") "6.88%"
Set-TextValue $ws.Range("E41") "5.71%"
Set-TextValue $ws.Range("D42") "0.1412"
Set-TextValue $ws.Range("E42") "6.42%"
Set-TextValue $ws.Range("D43") "0.007104"
Set-TextValue $ws.Range("E43") "-26.52%"
Set-TextValue $ws.Range("D44") "0.002141"
Set-TextValue $ws.Range("E44") "-7.81%"
Set-TextValue $ws.Range("D45") "0.009182"
Set-TextValue $ws.Range("E45") "0.45%"
Set-TextValue $ws.Range("D46") "0.00006123"
Set-TextValue $ws.Range("E46") "0.47%"
Set-TextValue $ws.Range("E47") "-0.06%"
Set-TextValue $ws.Range("E48") "38.83%"
Set-TextValue $ws.Range("D49") "0.002001"
Set-TextValue $ws.Range("E49") "-0.06%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "-0.06%"
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "-0.06%"
